$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 431/432, pushing the existing rows 431-443
# down to 433-445 (matches the dimension growing from T443 to T445).
$ws.Rows("431:432").Insert()

# Row 431: new Clementina / Primera record for Región de O'Higgins
$ws.Cells.Item(431, 1).Value = 4
$ws.Cells.Item(431, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(431, 3).Value = "Los Lagos"
$ws.Cells.Item(431, 4).Value = 45075
$ws.Cells.Item(431, 5).Value = 10
$ws.Cells.Item(431, 6).Value = "Fruta"
$ws.Cells.Item(431, 7).Value = 100102
$ws.Cells.Item(431, 8).Value = "Cítricos"
$ws.Cells.Item(431, 9).Value = 100102004
$ws.Cells.Item(431, 10).Value = "Mandarina"
$ws.Cells.Item(431, 11).Value = "Clementina"
$ws.Cells.Item(431, 12).Value = "Primera"
$ws.Cells.Item(431, 13).Value = 400
$ws.Cells.Item(431, 14).Value = 11000
$ws.Cells.Item(431, 15).Value = 12000
$ws.Cells.Item(431, 16).Value = 11500
$ws.Cells.Item(431, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(431, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(431, 19).Value = 1150
$ws.Cells.Item(431, 20).Value = 10

# Row 432: new Clementina / Segunda record for Región de O'Higgins
$ws.Cells.Item(432, 1).Value = 4
$ws.Cells.Item(432, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(432, 3).Value = "Los Lagos"
$ws.Cells.Item(432, 4).Value = 45075
$ws.Cells.Item(432, 5).Value = 10
$ws.Cells.Item(432, 6).Value = "Fruta"
$ws.Cells.Item(432, 7).Value = 100102
$ws.Cells.Item(432, 8).Value = "Cítricos"
$ws.Cells.Item(432, 9).Value = 100102004
$ws.Cells.Item(432, 10).Value = "Mandarina"
$ws.Cells.Item(432, 11).Value = "Clementina"
$ws.Cells.Item(432, 12).Value = "Segunda"
$ws.Cells.Item(432, 13).Value = 200
$ws.Cells.Item(432, 14).Value = 9000
$ws.Cells.Item(432, 15).Value = 9000
$ws.Cells.Item(432, 16).Value = 9000
$ws.Cells.Item(432, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(432, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(432, 19).Value = 900
$ws.Cells.Item(432, 20).Value = 10
